# Apply the described update:
# 1. Update the "Förändrad" (C column) date from 2024-01-15 (45306) to 2024-01-16 (45307)
#    for all existing data rows (rows 2-27).
# 2. Ensure row 27 has an explicit row height (matches the other data rows).
# 3. Add a new data row (row 28) for case "A 1668-2024".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update column C (Förändrad) for rows 2 through 27 to the new date serial value 45307.
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 3).Value = 45307
}

# 2. Give row 27 the same explicit row height as the surrounding rows.
$ws.Rows(27).RowHeight = 15

# 3. Append the new row with case "A 1668-2024" at row 28.
$newRow = 28

$ws.Cells.Item($newRow, 1).Value = "A 1668-2024"      # A: Beteckning
$ws.Cells.Item($newRow, 2).Value = 45306               # B: Datum
$ws.Cells.Item($newRow, 3).Value = 45307               # C: Förändrad
$ws.Cells.Item($newRow, 4).Value = "OKÄNT"            # D: Län
$ws.Cells.Item($newRow, 5).Value = "OKÄNT"            # E: Kommun
$ws.Cells.Item($newRow, 7).Value = 0.5                 # G: Area (ha)
$ws.Cells.Item($newRow, 8).Value = 0                   # H: Fridlysta
$ws.Cells.Item($newRow, 9).Value = 0                   # I: Signalarter
$ws.Cells.Item($newRow, 10).Value = 0                  # J: NT
$ws.Cells.Item($newRow, 11).Value = 0                  # K: VU
$ws.Cells.Item($newRow, 12).Value = 0                  # L: EN
$ws.Cells.Item($newRow, 13).Value = 0                  # M: CR
$ws.Cells.Item($newRow, 14).Value = 0                  # N: RE
$ws.Cells.Item($newRow, 15).Value = 0                  # O: Rödlistade
$ws.Cells.Item($newRow, 16).Value = 0                  # P: Hotade
$ws.Cells.Item($newRow, 17).Value = 0                  # Q: Alla arter

# Apply the same date number format as the other rows' B/C columns.
$ws.Cells.Item($newRow, 2).NumberFormat = $ws.Cells.Item($newRow - 1, 2).NumberFormat
$ws.Cells.Item($newRow, 3).NumberFormat = $ws.Cells.Item($newRow - 1, 3).NumberFormat

# Apply the same wrap-text style as column R used on previous rows, keeping it empty.
$ws.Cells.Item($newRow, 18).WrapText = $true
